# Re-point the three data tables (slides 14-16) from the old custom
# table style to the built-in table style {1D1F2745-6E61-46BF-84BB-9A07B685E9F6}.
#
# Table.Style is read-only in this object model (PowerPoint requires
# Table.ApplyStyle("{guid}") to change it), so we locate the table shape
# on each slide and call ApplyStyle with the target style id.

$p = $ppt.ActivePresentation
$newStyleId = "{1D1F2745-6E61-46BF-84BB-9A07B685E9F6}"

foreach ($slideIndex in 14,15,16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle($newStyleId)
        }
    }
}
